# Functional Simulation of Bus Routes
# - Flip the deployment "index" flag columns on Routes so the first stop
#   (mean 0 / std 0) moves from row 2 to row 3 for each of the
#   recharge / refill / route_1 / route_2 sections.
# - Add a "Demands" sheet describing route demand parameters.
# - Add a "Possible_Updates" sheet with a note about the change, and make
#   it the active tab (matching Excel's tabSelected behaviour after adding
#   sheets at the end).

$wb = $excel.ActiveWorkbook
$routes = $wb.Worksheets.Item(1)

# --- Flip the 0/1 "index" flag columns (A, D, G, J) on Routes ---
$flipCols = @("A", "D", "G", "J")
foreach ($col in $flipCols) {
    for ($row = 2; $row -le 8; $row++) {
        $cell = $routes.Range("$col$row")
        $val = $cell.Value2
        if ($val -eq 0) {
            $cell.Value = 1
        } elseif ($val -eq 1) {
            $cell.Value = 0
        }
    }
}

# Move the remembered selection on Routes from K5 to K4 (as in the saved file)
$routes.Range("K4").Select() | Out-Null

# --- Add the "Demands" sheet right after "Routes" ---
$demands = $wb.Worksheets.Add($null, $routes)
$demands.Name = "Demands"

$demands.Range("A1").Value = "Route"
$demands.Range("B1").Value = "a"
$demands.Range("C1").Value = "b"
$demands.Range("D1").Value = "c"
$demands.Range("E1").Value = "d"
$demands.Range("F1").Value = "charge"

$demandRows = @(
    @(1, 2, 3, 0, 75, 25),
    @(2, 2, 3, -180, 120, 20),
    @(3, 1, 3, 0, 45, 15)
)

$r = 2
foreach ($row in $demandRows) {
    $demands.Range("A$r").Value = $row[0]
    $demands.Range("A$r").NumberFormat = "@"
    $demands.Range("B$r").Value = $row[1]
    $demands.Range("C$r").Value = $row[2]
    $demands.Range("D$r").Value = $row[3]
    $demands.Range("E$r").Value = $row[4]
    $demands.Range("F$r").Value = $row[5]
    $r++
}

$demands.Range("F4").Select() | Out-Null

# --- Add the "Possible_Updates" sheet at the end, it becomes the active tab ---
$updates = $wb.Worksheets.Add($null, $demands)
$updates.Name = "Possible_Updates"

$updates.Range("A1").Value = "* Updating 1st stop of bus as deployment with mean 0 and std 0"
$updates.Range("A2").Select() | Out-Null
